# Apply 2024-04-03 violent crime data updates (column K = year 2024 year-to-date totals),
# plus a handful of upstream corrections to 2022 (column I) and 2016 (column C) values
# that shipped in the same data refresh, across the Citywide Totals, By Neighborhood, and
# per-neighborhood worksheets.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "Citywide Totals"; Cell = "K2"; Value = 1785 },
    @{ Sheet = "Citywide Totals"; Cell = "K3"; Value = 1699 },
    @{ Sheet = "Citywide Totals"; Cell = "C4"; Value = 1848 },
    @{ Sheet = "Citywide Totals"; Cell = "K4"; Value = 366 },
    @{ Sheet = "Citywide Totals"; Cell = "K5"; Value = 115 },
    @{ Sheet = "Citywide Totals"; Cell = "I6"; Value = 8962 },
    @{ Sheet = "Citywide Totals"; Cell = "K6"; Value = 2202 },
    @{ Sheet = "Citywide Totals"; Cell = "C7"; Value = 28392 },
    @{ Sheet = "Citywide Totals"; Cell = "I7"; Value = 26238 },
    @{ Sheet = "Citywide Totals"; Cell = "K7"; Value = 6167 },
    @{ Sheet = "Austin"; Cell = "K2"; Value = 119 },
    @{ Sheet = "Austin"; Cell = "K6"; Value = 141 },
    @{ Sheet = "Austin"; Cell = "K7"; Value = 407 },
    @{ Sheet = "South Chicago"; Cell = "K2"; Value = 53 },
    @{ Sheet = "South Chicago"; Cell = "K7"; Value = 129 },
    @{ Sheet = "Garfield Park"; Cell = "K3"; Value = 95 },
    @{ Sheet = "Garfield Park"; Cell = "K5"; Value = 5 },
    @{ Sheet = "Garfield Park"; Cell = "K6"; Value = 63 },
    @{ Sheet = "Garfield Park"; Cell = "K7"; Value = 252 },
    @{ Sheet = "West Pullman"; Cell = "K3"; Value = 31 },
    @{ Sheet = "West Pullman"; Cell = "K7"; Value = 106 },
    @{ Sheet = "Grand Crossing"; Cell = "K6"; Value = 73 },
    @{ Sheet = "Grand Crossing"; Cell = "K7"; Value = 207 },
    @{ Sheet = "New City"; Cell = "K3"; Value = 39 },
    @{ Sheet = "New City"; Cell = "K7"; Value = 152 },
    @{ Sheet = "By Neighborhood"; Cell = "K6"; Value = 50 },
    @{ Sheet = "By Neighborhood"; Cell = "K7"; Value = 173 },
    @{ Sheet = "By Neighborhood"; Cell = "K8"; Value = 407 },
    @{ Sheet = "By Neighborhood"; Cell = "K11"; Value = 132 },
    @{ Sheet = "By Neighborhood"; Cell = "K16"; Value = 13 },
    @{ Sheet = "By Neighborhood"; Cell = "K19"; Value = 170 },
    @{ Sheet = "By Neighborhood"; Cell = "K20"; Value = 136 },
    @{ Sheet = "By Neighborhood"; Cell = "K21"; Value = 19 },
    @{ Sheet = "By Neighborhood"; Cell = "K26"; Value = 9 },
    @{ Sheet = "By Neighborhood"; Cell = "K29"; Value = 299 },
    @{ Sheet = "By Neighborhood"; Cell = "K31"; Value = 70 },
    @{ Sheet = "By Neighborhood"; Cell = "K33"; Value = 252 },
    @{ Sheet = "By Neighborhood"; Cell = "K35"; Value = 10 },
    @{ Sheet = "By Neighborhood"; Cell = "K37"; Value = 207 },
    @{ Sheet = "By Neighborhood"; Cell = "K42"; Value = 211 },
    @{ Sheet = "By Neighborhood"; Cell = "K44"; Value = 58 },
    @{ Sheet = "By Neighborhood"; Cell = "K48"; Value = 68 },
    @{ Sheet = "By Neighborhood"; Cell = "K49"; Value = 37 },
    @{ Sheet = "By Neighborhood"; Cell = "K52"; Value = 169 },
    @{ Sheet = "By Neighborhood"; Cell = "K60"; Value = 43 },
    @{ Sheet = "By Neighborhood"; Cell = "C63"; Value = 276 },
    @{ Sheet = "By Neighborhood"; Cell = "K63"; Value = 18 },
    @{ Sheet = "By Neighborhood"; Cell = "K65"; Value = 152 },
    @{ Sheet = "By Neighborhood"; Cell = "K67"; Value = 240 },
    @{ Sheet = "By Neighborhood"; Cell = "K73"; Value = 61 },
    @{ Sheet = "By Neighborhood"; Cell = "K77"; Value = 41 },
    @{ Sheet = "By Neighborhood"; Cell = "K78"; Value = 79 },
    @{ Sheet = "By Neighborhood"; Cell = "K79"; Value = 164 },
    @{ Sheet = "By Neighborhood"; Cell = "K83"; Value = 129 },
    @{ Sheet = "By Neighborhood"; Cell = "K85"; Value = 306 },
    @{ Sheet = "By Neighborhood"; Cell = "K88"; Value = 79 },
    @{ Sheet = "By Neighborhood"; Cell = "K89"; Value = 79 },
    @{ Sheet = "By Neighborhood"; Cell = "K90"; Value = 53 },
    @{ Sheet = "By Neighborhood"; Cell = "K91"; Value = 54 },
    @{ Sheet = "By Neighborhood"; Cell = "K95"; Value = 106 },
    @{ Sheet = "By Neighborhood"; Cell = "I97"; Value = 236 },
    @{ Sheet = "By Neighborhood"; Cell = "K97"; Value = 53 },
    @{ Sheet = "By Neighborhood"; Cell = "K98"; Value = 43 },
    @{ Sheet = "By Neighborhood"; Cell = "C101"; Value = 28392 },
    @{ Sheet = "By Neighborhood"; Cell = "I101"; Value = 26238 },
    @{ Sheet = "By Neighborhood"; Cell = "K101"; Value = 6167 },
    @{ Sheet = "Gage Park"; Cell = "K6"; Value = 29 },
    @{ Sheet = "Gage Park"; Cell = "K7"; Value = 70 },
    @{ Sheet = "North Lawndale"; Cell = "K4"; Value = 14 },
    @{ Sheet = "North Lawndale"; Cell = "K6"; Value = 81 },
    @{ Sheet = "North Lawndale"; Cell = "K7"; Value = 240 },
    @{ Sheet = "Lincoln Park"; Cell = "K2"; Value = 3 },
    @{ Sheet = "Lincoln Park"; Cell = "K7"; Value = 37 },
    @{ Sheet = "Englewood"; Cell = "K2"; Value = 81 },
    @{ Sheet = "Englewood"; Cell = "K3"; Value = 101 },
    @{ Sheet = "Englewood"; Cell = "K7"; Value = 299 },
    @{ Sheet = "Lake View"; Cell = "K6"; Value = 31 },
    @{ Sheet = "Lake View"; Cell = "K7"; Value = 68 },
    @{ Sheet = "Chatham"; Cell = "K2"; Value = 54 },
    @{ Sheet = "Chatham"; Cell = "K4"; Value = 3 },
    @{ Sheet = "Chatham"; Cell = "K7"; Value = 170 },
    @{ Sheet = "Irving Park"; Cell = "K3"; Value = 19 },
    @{ Sheet = "Irving Park"; Cell = "K4"; Value = 4 },
    @{ Sheet = "Irving Park"; Cell = "K7"; Value = 58 },
    @{ Sheet = "Ashburn"; Cell = "K6"; Value = 18 },
    @{ Sheet = "Ashburn"; Cell = "K7"; Value = 50 },
    @{ Sheet = "Humboldt Park"; Cell = "K6"; Value = 94 },
    @{ Sheet = "Humboldt Park"; Cell = "K7"; Value = 211 },
    @{ Sheet = "Rogers Park"; Cell = "K2"; Value = 25 },
    @{ Sheet = "Rogers Park"; Cell = "K6"; Value = 28 },
    @{ Sheet = "Rogers Park"; Cell = "K7"; Value = 79 },
    @{ Sheet = "Washington Park"; Cell = "K2"; Value = 18 },
    @{ Sheet = "Washington Park"; Cell = "K7"; Value = 54 },
    @{ Sheet = "Chinatown"; Cell = "K3"; Value = 6 },
    @{ Sheet = "Chinatown"; Cell = "K7"; Value = 19 },
    @{ Sheet = "Roseland"; Cell = "K2"; Value = 58 },
    @{ Sheet = "Roseland"; Cell = "K7"; Value = 164 },
    @{ Sheet = "Chicago Lawn"; Cell = "K2"; Value = 39 },
    @{ Sheet = "Chicago Lawn"; Cell = "K6"; Value = 50 },
    @{ Sheet = "Chicago Lawn"; Cell = "K7"; Value = 136 },
    @{ Sheet = "Auburn Gresham"; Cell = "K4"; Value = 7 },
    @{ Sheet = "Auburn Gresham"; Cell = "K6"; Value = 44 },
    @{ Sheet = "Auburn Gresham"; Cell = "K7"; Value = 173 },
    @{ Sheet = "Wicker Park"; Cell = "K4"; Value = 3 },
    @{ Sheet = "Wicker Park"; Cell = "K6"; Value = 32 },
    @{ Sheet = "Wicker Park"; Cell = "K7"; Value = 43 },
    @{ Sheet = "East Village"; Cell = "K6"; Value = 6 },
    @{ Sheet = "East Village"; Cell = "K7"; Value = 9 },
    @{ Sheet = "Belmont Cragin"; Cell = "K6"; Value = 53 },
    @{ Sheet = "Belmont Cragin"; Cell = "K7"; Value = 132 },
    @{ Sheet = "Gold Coast"; Cell = "K3"; Value = 2 },
    @{ Sheet = "Gold Coast"; Cell = "K7"; Value = 10 },
    @{ Sheet = "Portage Park"; Cell = "K2"; Value = 17 },
    @{ Sheet = "Portage Park"; Cell = "K7"; Value = 61 },
    @{ Sheet = "West Town"; Cell = "I6"; Value = 156 },
    @{ Sheet = "West Town"; Cell = "K6"; Value = 35 },
    @{ Sheet = "West Town"; Cell = "I7"; Value = 236 },
    @{ Sheet = "West Town"; Cell = "K7"; Value = 53 },
    @{ Sheet = "United Center"; Cell = "K2"; Value = 18 },
    @{ Sheet = "United Center"; Cell = "K6"; Value = 45 },
    @{ Sheet = "United Center"; Cell = "K7"; Value = 79 },
    @{ Sheet = "Uptown"; Cell = "K3"; Value = 28 },
    @{ Sheet = "Uptown"; Cell = "K7"; Value = 79 },
    @{ Sheet = "Washington Heights"; Cell = "K2"; Value = 24 },
    @{ Sheet = "Washington Heights"; Cell = "K7"; Value = 53 },
    @{ Sheet = "Morgan Park"; Cell = "K2"; Value = 9 },
    @{ Sheet = "Morgan Park"; Cell = "K7"; Value = 43 },
    @{ Sheet = "South Shore"; Cell = "K2"; Value = 111 },
    @{ Sheet = "South Shore"; Cell = "K3"; Value = 100 },
    @{ Sheet = "South Shore"; Cell = "K7"; Value = 306 },
    @{ Sheet = "Riverdale"; Cell = "K2"; Value = 20 },
    @{ Sheet = "Riverdale"; Cell = "K7"; Value = 41 },
    @{ Sheet = "Little Village"; Cell = "K2"; Value = 40 },
    @{ Sheet = "Little Village"; Cell = "K3"; Value = 37 },
    @{ Sheet = "Little Village"; Cell = "K6"; Value = 78 },
    @{ Sheet = "Little Village"; Cell = "K7"; Value = 169 },
    @{ Sheet = "Bucktown"; Cell = "K6"; Value = 10 },
    @{ Sheet = "Bucktown"; Cell = "K7"; Value = 13 }
)

$seenSheets = @{}
foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
    $seenSheets[$u.Sheet] = $true
}

Write-Host ("Applied {0} cell updates across {1} worksheets." -f $updates.Count, $seenSheets.Keys.Count)
